# Updated cryptos list on Wed Feb 28 11:35:28 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.368.59'
$ws.Range("E2").Value = '  +4.99%  '
$ws.Range("D3").Value = '3.348.46'
$ws.Range("E3").Value = '  +2.84%  '
$ws.Range("E4").Value = '  +0.01%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = '''410.94'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +2.84%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = '''113.88'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +2.39%  '
$style = $ws.Range("D7").Style
$ws.Range("D7").Value = '''0.591'
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  +5.81%  '
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  +0.02%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = '''0.641'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  +3.51%  '
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = '''40.43'
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  +2.23%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = '''0.0992'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  +4.81%  '
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = '''0.144'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("D13").Value = '3.878.20'
$ws.Range("E13").Value = '  +2.88%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = '''8.58'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +5.64%  '
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = '''19.58'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").Value = '3.340.45'
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '59.179.95'
$ws.Range("E18").Value = '  +4.86%  '
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = '''10.77'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  -2.11%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").Value = '''3.39'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("E21").Value = '  +7.58%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = '''13.36'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +2.33%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").Value = '''306.71'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("E24").Value = '  +1.03%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = '''3.22'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -0.30%  '
$style = $ws.Range("D26").Style
$ws.Range("D26").Value = '''28.75'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  +1.82%  '
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = '''7.76'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  +5.55%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = '''0.181'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  +6.63%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = '''7.97'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -1.75%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").Value = '''0.117'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  +4.32%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = '''40.23'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +5.18%  '
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = '''0.0519'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  +6.12%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = '''52.26'
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  +1.05%  '
$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$style = $ws.Range("D36").Style
$ws.Range("D36").Value = '''2.11'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  -1.65%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = '''3.15'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +0.43%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = '''0.999'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -2.34%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = '''137.99'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  +3.34%  '
$ws.Range("E41").Value = '  +2.58%  '
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = '''17.21'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = '''3.97'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  -0.52%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = '''0.286'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +0.27%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = '''22.67'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +1.91%  '
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = '''2.26'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  +8.30%  '
$ws.Range("D48").Value = '2.225.95'
$ws.Range("E48").Value = '  +3.52%  '
$ws.Range("E49").Value = '  -0.84%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = '''1.92'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -9.47%  '
$ws.Range("E51").Value = '  +6.21%  '
